$wb = $excel.ActiveWorkbook

# This script applies numeric value corrections to several rows across
# multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), as produced by
# the scheduled market-data refresh run.

# --- ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H87").Value = 18678.46
$ws.Range("J87").Value = 18678.46
$ws.Range("L87").Value = 18678.46
$ws.Range("N87").Value = -21174.46
$ws.Range("H90").Value = 18678.46
$ws.Range("J90").Value = 18678.46
$ws.Range("L90").Value = 56035.38
$ws.Range("N90").Value = -68515.38
$ws.Range("H112").Value = 1434.6
$ws.Range("I112").Value = 423.16666
$ws.Range("J112").Value = 1754
$ws.Range("K112").Value = 1269.49998
$ws.Range("L112").Value = 5262
$ws.Range("M112").Value = -161.4999800000001
$ws.Range("N112").Value = -7478
$ws.Range("H125").Value = 21900
$ws.Range("I125").Value = 100000
$ws.Range("J125").Value = 2375
$ws.Range("K125").Value = 900000
$ws.Range("L125").Value = 21375
$ws.Range("M125").Value = -897540
$ws.Range("N125").Value = -26295
$ws.Range("H131").Value = 9434.833000000001
$ws.Range("I131").Value = 10028.909
$ws.Range("J131").Value = 2900
$ws.Range("K131").Value = 30086.727
$ws.Range("L131").Value = 8700
$ws.Range("M131").Value = -25046.727
$ws.Range("N131").Value = -18780

# --- ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 5523.7617
$ws.Range("I32").Value = 4058.7942
$ws.Range("J32").Value = 11749.875
$ws.Range("K32").Value = 4058.7942
$ws.Range("L32").Value = 11749.875
$ws.Range("M32").Value = -3771.7942
$ws.Range("N32").Value = -12323.875
$ws.Range("H45").Value = 1751
$ws.Range("I45").Value = 1834
$ws.Range("J45").Value = 1557.3334
$ws.Range("K45").Value = 1834
$ws.Range("L45").Value = 1557.3334
$ws.Range("M45").Value = -1457
$ws.Range("N45").Value = -2311.3334
$ws.Range("H61").Value = 4974.857
$ws.Range("I61").Value = 4974.857
$ws.Range("K61").Value = 4974.857
$ws.Range("M61").Value = -4762.857
$ws.Range("H136").Value = 4974.857
$ws.Range("I136").Value = 4974.857
$ws.Range("K136").Value = 14924.571
$ws.Range("M136").Value = -12374.571

# --- BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H105").Value = 23811888
$ws.Range("I105").Value = 29414196
$ws.Range("J105").Value = 2075.125
$ws.Range("K105").Value = 29414196
$ws.Range("L105").Value = 2075.125
$ws.Range("M105").Value = -29412449
$ws.Range("N105").Value = -5569.125
$ws.Range("H107").Value = 2159.5386
$ws.Range("I107").Value = 2097.2856
$ws.Range("J107").Value = 2232.1667
$ws.Range("K107").Value = 2097.2856
$ws.Range("L107").Value = 2232.1667
$ws.Range("M107").Value = -177.2856000000002
$ws.Range("N107").Value = -6072.1667

# --- CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 39643.59
$ws.Range("I31").Value = 44300.734
$ws.Range("J31").Value = 4715
$ws.Range("K31").Value = 44300.734
$ws.Range("L31").Value = 4715
$ws.Range("M31").Value = -44005.734
$ws.Range("N31").Value = -5305
$ws.Range("H34").Value = 39643.59
$ws.Range("I34").Value = 44300.734
$ws.Range("J34").Value = 4715
$ws.Range("K34").Value = 44300.734
$ws.Range("L34").Value = 4715
$ws.Range("M34").Value = -44098.734
$ws.Range("N34").Value = -5119

# --- CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H2").Value = 266.7619
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 549.2
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 3295.2
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -3521.2
$ws.Range("H92").Value = 158.51724
$ws.Range("I92").Value = 166.7619
$ws.Range("J92").Value = 136.875
$ws.Range("K92").Value = 500.2857
$ws.Range("L92").Value = 410.625
$ws.Range("M92").Value = 747.7143
$ws.Range("N92").Value = -2906.625
$ws.Range("H97").Value = 1368.5652
$ws.Range("I97").Value = 993.8
$ws.Range("J97").Value = 1472.6666
$ws.Range("K97").Value = 2981.4
$ws.Range("L97").Value = 4417.9998
$ws.Range("M97").Value = -2485.4
$ws.Range("N97").Value = -5409.9998
$ws.Range("H98").Value = 667232.7
$ws.Range("I98").Value = 402.55554
$ws.Range("J98").Value = 1667477.9
$ws.Range("K98").Value = 1207.66662
$ws.Range("L98").Value = 5002433.699999999
$ws.Range("M98").Value = 290.33338
$ws.Range("N98").Value = -5005429.699999999
$ws.Range("H107").Value = 8642.691999999999
$ws.Range("I107").Value = 25565
$ws.Range("J107").Value = 1121.6666
$ws.Range("K107").Value = 76695
$ws.Range("L107").Value = 3364.9998
$ws.Range("M107").Value = -74775
$ws.Range("N107").Value = -7204.9998

# --- GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 2583.3333
$ws.Range("I80").Value = 2625
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2625
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1627
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 2583.3333
$ws.Range("I83").Value = 2625
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 13125
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -8133
$ws.Range("N83").Value = -22484
$ws.Range("H97").Value = 2329.4443
$ws.Range("I97").Value = 1339.909
$ws.Range("J97").Value = 3884.4285
$ws.Range("K97").Value = 1339.909
$ws.Range("L97").Value = 3884.4285
$ws.Range("M97").Value = -843.9090000000001
$ws.Range("N97").Value = -4876.4285

# --- LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H132").Value = 7101.25
$ws.Range("I132").Value = 8810.143
$ws.Range("J132").Value = 4708.8
$ws.Range("K132").Value = 26430.429
$ws.Range("L132").Value = 14126.4
$ws.Range("M132").Value = -23900.429
$ws.Range("N132").Value = -19186.4
$ws.Range("H133").Value = 32947.332
$ws.Range("J133").Value = 32947.332
$ws.Range("L133").Value = 32947.332
$ws.Range("N133").Value = -38007.332

# --- WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H113").Value = 45000180
$ws.Range("I113").Value = 45454748
$ws.Range("J113").Value = 41666690
$ws.Range("K113").Value = 136364244
$ws.Range("L113").Value = 125000070
$ws.Range("M113").Value = -136362074
$ws.Range("N113").Value = -125004410
$ws.Range("H136").Value = 10521.454
$ws.Range("I136").Value = 10521.454
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 31564.362
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -29014.362
$ws.Range("N136").ClearContents()
